# Updates cryptos list values (Price and Volume(1h) columns) per commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be written as literal text so Excel does not
    # "smart" convert numeric-looking strings (e.g. "88.40" -> 88.4)
    # or percent-looking strings, then drop the explicit format again
    # so the cell keeps no style override (matches original file).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "40.066.09"
Set-TextValue $ws.Range("E2") "  +0.45%  "
Set-TextValue $ws.Range("D3") "2.216.49"
Set-TextValue $ws.Range("E3") "  -0.68%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "290.89"
Set-TextValue $ws.Range("E5") "  -2.62%  "
Set-TextValue $ws.Range("D6") "88.40"
Set-TextValue $ws.Range("E6") "  +4.68%  "
Set-TextValue $ws.Range("E7") "  +0.21%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.472"
Set-TextValue $ws.Range("E9") "  +0.78%  "
Set-TextValue $ws.Range("D10") "30.84"
Set-TextValue $ws.Range("E10") "  +3.05%  "
Set-TextValue $ws.Range("E11") "  +0.44%  "
Set-TextValue $ws.Range("D12") "47.77"
Set-TextValue $ws.Range("E12") "  +2.47%  "
Set-TextValue $ws.Range("D13") "0.110"
Set-TextValue $ws.Range("E13") "  +2.11%  "
Set-TextValue $ws.Range("D14") "6.45"
Set-TextValue $ws.Range("E14") "  +2.06%  "
Set-TextValue $ws.Range("D15") "2.558.91"
Set-TextValue $ws.Range("E15") "  -0.53%  "
Set-TextValue $ws.Range("D16") "14.02"
Set-TextValue $ws.Range("E16") "  -1.00%  "
Set-TextValue $ws.Range("D17") "2.211.14"
Set-TextValue $ws.Range("E17") "  -0.59%  "
Set-TextValue $ws.Range("D18") "0.730"
Set-TextValue $ws.Range("E18") "  +1.42%  "
Set-TextValue $ws.Range("D19") "39.983.77"
Set-TextValue $ws.Range("E19") "  +0.53%  "
Set-TextValue $ws.Range("D20") "11.76"
Set-TextValue $ws.Range("E20") "  +11.51%  "
Set-TextValue $ws.Range("D21") "0.0₃0887"
Set-TextValue $ws.Range("E21") "  +0.71%  "
Set-TextValue $ws.Range("D22") "5.82"
Set-TextValue $ws.Range("E22") "  +0.57%  "
Set-TextValue $ws.Range("D23") "65.65"
Set-TextValue $ws.Range("E23") "  +0.64%  "
Set-TextValue $ws.Range("D24") "236.94"
Set-TextValue $ws.Range("E24") "  +1.04%  "
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  -0.12%  "
Set-TextValue $ws.Range("E26") "  +2.01%  "
Set-TextValue $ws.Range("D27") "1.85"
Set-TextValue $ws.Range("E27") "  +1.48%  "
Set-TextValue $ws.Range("D28") "22.65"
Set-TextValue $ws.Range("E28") "  -0.71%  "
Set-TextValue $ws.Range("D29") "2.20"
Set-TextValue $ws.Range("E29") "  +4.31%  "
Set-TextValue $ws.Range("E30") "  +0.63%  "
Set-TextValue $ws.Range("D31") "153.42"
Set-TextValue $ws.Range("E31") "  +2.34%  "
Set-TextValue $ws.Range("D32") "32.29"
Set-TextValue $ws.Range("E32") "  -0.87%  "
Set-TextValue $ws.Range("D34") "4.97"
Set-TextValue $ws.Range("E34") "  +2.43%  "
Set-TextValue $ws.Range("D35") "0.0720"
Set-TextValue $ws.Range("E35") "  +2.39%  "
Set-TextValue $ws.Range("E36") "  -0.71%  "
Set-TextValue $ws.Range("E37") "  +6.70%  "
Set-TextValue $ws.Range("D38") "16.22"
Set-TextValue $ws.Range("E38") "  -1.82%  "
Set-TextValue $ws.Range("E39") "  +0.61%  "
Set-TextValue $ws.Range("D40") "0.100"
Set-TextValue $ws.Range("E40") "  +2.29%  "
Set-TextValue $ws.Range("E41") "  +2.84%  "
Set-TextValue $ws.Range("D42") "2.093.28"
Set-TextValue $ws.Range("E42") "  +8.20%  "
Set-TextValue $ws.Range("D43") "3.84"
Set-TextValue $ws.Range("E43") "  +4.62%  "
Set-TextValue $ws.Range("E44") "  +1.33%  "
Set-TextValue $ws.Range("E45") "  +1.59%  "
Set-TextValue $ws.Range("D46") "17.88"
Set-TextValue $ws.Range("E46") "  +8.20%  "
Set-TextValue $ws.Range("D47") "9.94"
Set-TextValue $ws.Range("E47") "  +7.46%  "
Set-TextValue $ws.Range("D48") "2.67"
Set-TextValue $ws.Range("E48") "  +2.04%  "
Set-TextValue $ws.Range("D49") "2.432.38"
Set-TextValue $ws.Range("E49") "  -0.45%  "
Set-TextValue $ws.Range("D50") "69.91"
Set-TextValue $ws.Range("E50") "  -1.66%  "
Set-TextValue $ws.Range("E51") "  +4.44%  "
